$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 - this shifts existing rows 8..35 down to 9..36
# and automatically extends the used range / dimension to A1:R36.
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with the new record's data.
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C8").Value = "Los Lagos"
$ws.Range("D8").Value = 44624
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 100112030
$ws.Range("G8").Value = "Poroto granado"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 33000
$ws.Range("L8").Value = 33000
$ws.Range("M8").Value = 33000
$ws.Range("N8").Value = '$/saco 25 kilos'
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 1320
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
